# Excel COM-interop script: "funnel fix scale fig"
# - Corrects the "subgroup" (Moderator) classification in column I for several studies
# - Adds a new helper column J computing (N2 - N1), i.e. F - C, the sample-size delta
#   used to scale the funnel-plot figure
# - Updates the active selection to reflect where the user ended up working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the "subgroup" (Moderator) values in column I ----------------------
$ws.Range("I2").Value  = "subgroup1"
$ws.Range("I3").Value  = "subgroup1"
$ws.Range("I4").Value  = "subgroup1"
$ws.Range("I5").Value  = "subgroup2"
$ws.Range("I6").Value  = "subgroup2"
$ws.Range("I7").Value  = "subgroup3"
$ws.Range("I8").Value  = "subgroup1"
$ws.Range("I9").Value  = "subgroup1"
$ws.Range("I10").Value = "subgroup1"
$ws.Range("I11").Value = "subgroup2"

# --- Add new column J with the N2-N1 scale-fix formula -------------------------
# J2 is entered on its own; J3:J11 are filled together afterwards so Excel records
# them as one shared formula group (mirrors a fill-down of the first formula).
$ws.Range("J2").Formula = "=F2-C2"
$ws.Range("J3:J11").Formula = "=F3-C3"

# --- Update view state: scroll position / selection -----------------------------
# Scroll the window so column B is the left-most visible column, then leave the
# selection on I4 (matching where editing finished).
$ws.Activate()
$null = $ws.Range("B1").Select()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollColumn = 2
    $win.ScrollRow = 1
}
$null = $ws.Range("I4").Select()
